$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings (some look numeric, e.g. "1.003" or "313.90")
# Force the column to Text format first so Excel does not coerce these into numbers,
# matching the original workbook where every D-column cell is a text (inlineStr) value.
$ws.Range('D2:D51').NumberFormat = "@"

$ws.Range('D2').Value = '27.489.77'
$ws.Range('E2').Value = '  -0.89%  '
$ws.Range('D3').Value = '1.830.77'
$ws.Range('E3').Value = '  -0.92%  '
$ws.Range('D4').Value = '1.003'
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').Value = '313.90'
$ws.Range('E5').Value = '  +0.00%  '
$ws.Range('E6').Value = '  +0.04%  '
$ws.Range('D7').Value = '0.4289'
$ws.Range('E7').Value = '  -0.95%  '
$ws.Range('D8').Value = '0.3655'
$ws.Range('E8').Value = '  +0.10%  '
$ws.Range('D9').Value = '0.07274'
$ws.Range('E9').Value = '  -0.83%  '
$ws.Range('D10').Value = '0.8678'
$ws.Range('E10').Value = '  -0.96%  '
$ws.Range('E11').Value = '  -0.61%  '
$ws.Range('D12').Value = '1.845.97'
$ws.Range('E12').Value = '  +1.10%  '
$ws.Range('E13').Value = '  +1.22%  '
$ws.Range('D14').Value = '6.530'
$ws.Range('E14').Value = '  +0.19%  '
$ws.Range('D15').Value = '0.06938'
$ws.Range('E15').Value = '  +0.22%  '
$ws.Range('D16').Value = '1.004'
$ws.Range('E16').Value = '  +0.23%  '
$ws.Range('D17').Value = '80.52'
$ws.Range('E17').Value = '  +0.72%  '
$ws.Range('D18').Value = '0.000008890'
$ws.Range('E18').Value = '  -1.09%  '
$ws.Range('D19').Value = '1.003'
$ws.Range('E19').Value = '  +0.26%  '
$ws.Range('E20').Value = '  +0.17%  '
$ws.Range('D21').Value = '27.844.57'
$ws.Range('E21').Value = '  +0.86%  '
$ws.Range('D22').Value = '5.137'
$ws.Range('E22').Value = '  +3.20%  '
$ws.Range('E23').Value = '  +4.00%  '
$ws.Range('D24').Value = '2.105.62'
$ws.Range('E24').Value = '  +2.77%  '
$ws.Range('D25').Value = '1.981'
$ws.Range('E25').Value = '  +0.04%  '
$ws.Range('D26').Value = '154.28'
$ws.Range('E26').Value = '  -1.15%  '
$ws.Range('D27').Value = '18.80'
$ws.Range('E27').Value = '  +0.97%  '
$ws.Range('D28').Value = '5.143'
$ws.Range('E28').Value = '  -2.03%  '
$ws.Range('E29').Value = '  -5.26%  '
$ws.Range('D30').Value = '1.829'
$ws.Range('E30').Value = '  -1.37%  '
$ws.Range('D31').Value = '0.08877'
$ws.Range('E31').Value = '  -0.29%  '
$ws.Range('D32').Value = '0.7539'
$ws.Range('E32').Value = '  +0.35%  '
$ws.Range('D33').Value = '2.983'
$ws.Range('E33').Value = '  +0.62%  '
$ws.Range('D34').Value = '4.538'
$ws.Range('E34').Value = '  +0.00%  '
$ws.Range('E35').Value = '  +1.11%  '
$ws.Range('E36').Value = '  +0.08%  '
$ws.Range('D37').Value = '1.094'
$ws.Range('D38').Value = '0.05308'
$ws.Range('E38').Value = '  -1.81%  '
$ws.Range('D39').Value = '0.01936'
$ws.Range('E39').Value = '  +0.31%  '
$ws.Range('D41').Value = '0.1663'
$ws.Range('E41').Value = '  +0.46%  '
$ws.Range('D42').Value = '0.5060'
$ws.Range('E42').Value = '  -0.47%  '
$ws.Range('D43').Value = '6.599'
$ws.Range('E43').Value = '  -0.87%  '
$ws.Range('D44').Value = '8.374'
$ws.Range('E44').Value = '  +0.66%  '
$ws.Range('E45').Value = '  +1.93%  '
$ws.Range('D46').Value = '105.93'
$ws.Range('E46').Value = '  +1.57%  '
$ws.Range('D47').Value = '0.06497'
$ws.Range('E47').Value = '  -0.61%  '
$ws.Range('D48').Value = '0.4677'
$ws.Range('E48').Value = '  +0.33%  '
$ws.Range('E49').Value = '  +0.03%  '
$ws.Range('D50').Value = '1.607'
$ws.Range('E50').Value = '  -0.99%  '
$ws.Range('D51').Value = '64.04'
$ws.Range('E51').Value = '  -0.37%  '
